{"js": "// Only check reasonable fear box and provide explanation if respondent not\n// convicted of SA.\n//\n// The template paragraph currently reads:\n//   {%p if potential_sexual_assault_exp %}\n// and must become:\n//   {%p if not respondent_sexual_assault_conviction %}\n//\n// (A second, unrelated occurrence of \"potential_sexual_assault_exp\" further\n// down, inside \"{{ potential_sexual_assault_exp }}\", must stay untouched.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.trim() === \"{%p if potential_sexual_assault_exp %}\"\n);\n\nif (!target) {\n  throw new Error(\"Could not locate target paragraph '{%p if potential_sexual_assault_exp %}'\");\n}\n\n// Search only within that paragraph so the later \"{{ potential_sexual_assault_exp }}\"\n// usage is left alone.\nconst matches = target.search(\"potential_sexual_assault_exp\", { matchCase: true });\nmatches.load(\"text\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not locate 'potential_sexual_assault_exp' inside target paragraph\");\n}\n\nmatches.items[0].insertText(\"not respondent_sexual_assault_conviction\", \"Replace\");\nawait context.sync();\n", "ps1": "# Only check reasonable fear box and provide explanation if respondent not\n# convicted of SA.\n#\n# The template paragraph currently reads:\n#   {%p if potential_sexual_assault_exp %}\n# and must become:\n#   {%p if not respondent_sexual_assault_conviction %}\n#\n# (A second, unrelated occurrence of \"potential_sexual_assault_exp\" further\n# down, inside \"{{ potential_sexual_assault_exp }}\", must stay untouched, so\n# the Find/Replace is scoped to just the target paragraph's range.)\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"{%p if potential_sexual_assault_exp %}\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate target paragraph '{%p if potential_sexual_assault_exp %}'\"\n}\n\n$r = $target.Range\n$find = $r.Find\n$find.ClearFormatting()\n$find.Text = \"potential_sexual_assault_exp\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"not respondent_sexual_assault_conviction\"\n\n# wdFindContinue = 1, wdReplaceOne = 1\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n\n$d.Save()\n"}
